$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.204710006713867
$ws.Range("B1").Value = 3.37756872177124
$ws.Range("C1").Value = 2.966285467147827
$ws.Range("D1").Value = 3.222946882247925
$ws.Range("E1").Value = 2.110574960708618
